$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_IB전략컨설팅부")

# Delete the four rows (큐로셀, 메가터치, 컨텍, 비아이매트릭스) that were removed
# from the strategy table. Delete bottom-up so row numbers stay valid.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(18).Delete()
